# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) with the latest scraped figures.
# Price cells that would otherwise be auto-parsed as numbers are entered
# with a leading apostrophe (Excel "text" quote-prefix) and then restored
# to the Normal cell style so the stored value stays plain text, matching
# the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.428.68'
$ws.Range('E2').Value = '  -2.43%  '
$ws.Range('D3').Value = '3.692.86'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'692.13"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = "'162.16"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.59%  '
$ws.Range('D7').Value = '3.691.48'
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -4.69%  '
$ws.Range('E10').Value = '  -8.20%  '
$ws.Range('D11').Value = "'7.39"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').Value = "'0.443"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.21%  '
$ws.Range('D13').Value = "'0.0000239"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.18%  '
$ws.Range('D14').Value = "'33.33"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.19%  '
$ws.Range('D15').Value = '4.314.38'
$ws.Range('E15').Value = '  -3.09%  '
$ws.Range('D16').Value = '3.694.87'
$ws.Range('E16').Value = '  -3.90%  '
$ws.Range('D17').Value = '69.463.63'
$ws.Range('E17').Value = '  -2.42%  '
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = "'16.16"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.68%  '
$ws.Range('D20').Value = "'6.57"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.98%  '
$ws.Range('D21').Value = "'479.88"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'10.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.78%  '
$ws.Range('E23').Value = '  -7.20%  '
$ws.Range('D24').Value = "'79.88"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.81%  '
$ws.Range('D25').Value = '3.837.79'
$ws.Range('D26').Value = "'0.0000130"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.09%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = "'11.38"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.79%  '
$ws.Range('D29').Value = "'9.47"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.59%  '
$ws.Range('E30').Value = '  -11.29%  '
$ws.Range('E31').Value = '  -10.16%  '
$ws.Range('D32').Value = "'6.84"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.54%  '
$ws.Range('E33').Value = '  -7.66%  '
$ws.Range('E34').Value = '  -5.40%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('E36').Value = '  -7.15%  '
$ws.Range('D37').Value = '3.655.60'
$ws.Range('E37').Value = '  -3.08%  '
$ws.Range('D38').Value = "'8.45"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.35%  '
$ws.Range('D39').Value = "'6.33"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.31%  '
$ws.Range('D40').Value = "'2.34"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('E41').Value = '  -8.26%  '
$ws.Range('E44').Value = '  -6.55%  '
$ws.Range('D45').Value = "'163.32"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.73%  '
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('D47').Value = "'30.13"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('E48').Value = '  -15.00%  '
$ws.Range('E49').Value = '  -2.53%  '
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('D51').Value = "'0.000283"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -9.02%  '
